$wb = $excel.ActiveWorkbook

# ALC row 61
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 1348
$ws.Range("I61").Value = 969.1429000000001
$ws.Range("K61").Value = 2907.4287
$ws.Range("M61").Value = -2735.4287

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3902.4546
$ws.Range("I86").Value = 4428.28
$ws.Range("J86").Value = 2259.25
$ws.Range("K86").Value = 4428.28
$ws.Range("L86").Value = 2259.25
$ws.Range("M86").Value = -3305.28
$ws.Range("N86").Value = -4505.25

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 3902.4546
$ws.Range("I89").Value = 4428.28
$ws.Range("J89").Value = 2259.25
$ws.Range("K89").Value = 22141.4
$ws.Range("L89").Value = 11296.25
$ws.Range("M89").Value = -16525.4
$ws.Range("N89").Value = -22528.25

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 3597.0588
$ws.Range("I107").Value = 1066.6154
$ws.Range("K107").Value = 1066.6154
$ws.Range("M107").Value = 853.3846000000001

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 10403.734
$ws.Range("I137").Value = 1772.4138
$ws.Range("J137").Value = 22919.15
$ws.Range("K137").Value = 5317.2414
$ws.Range("L137").Value = 68757.45000000001
$ws.Range("M137").Value = -2767.2414
$ws.Range("N137").Value = -73857.45000000001

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4995.909
$ws.Range("J141").Value = 4993.75
$ws.Range("L141").Value = 14981.25
$ws.Range("N141").Value = -25341.25

# ARM row 12
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 1073.5
$ws.Range("I12").Value = 1073.5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1073.5
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = $null
$ws.Range("N12").Value = -900.5

# ARM row 14
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 404.57144
$ws.Range("I14").Value = 398.66666
$ws.Range("J14").Value = 409
$ws.Range("K14").Value = 398.66666
$ws.Range("L14").Value = 409
$ws.Range("M14").Value = -223.66666
$ws.Range("N14").Value = -759

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4253.1455
$ws.Range("I32").Value = 1593.2307
$ws.Range("K32").Value = 1593.2307
$ws.Range("M32").Value = -1306.2307

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4544.273
$ws.Range("I102").Value = 4544.273
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4544.273
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = $null
$ws.Range("N102").Value = -2922.273

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 24740.105
$ws.Range("I20").Value = 9303.556
$ws.Range("J20").Value = 38633
$ws.Range("K20").Value = 9303.556
$ws.Range("L20").Value = 38633
$ws.Range("M20").Value = -9056.556
$ws.Range("N20").Value = -39127

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5002.25
$ws.Range("I86").Value = 2006
$ws.Range("K86").Value = 2006
$ws.Range("M86").Value = -883

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 5002.25
$ws.Range("I89").Value = 2006
$ws.Range("K89").Value = 10030
$ws.Range("M89").Value = -4414

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8508.076999999999
$ws.Range("I16").Value = 9012.916999999999
$ws.Range("K16").Value = 9012.916999999999
$ws.Range("M16").Value = -8725.916999999999

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 8508.076999999999
$ws.Range("I113").Value = 9012.916999999999
$ws.Range("K113").Value = 9012.916999999999
$ws.Range("M113").Value = -6842.916999999999

# CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 355.4
$ws.Range("I2").Value = 113
$ws.Range("J2").Value = 597.8
$ws.Range("K2").Value = 678
$ws.Range("L2").Value = 3586.8
$ws.Range("M2").Value = -565
$ws.Range("N2").Value = -3812.8

# CUL row 23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 425.15384
$ws.Range("I23").Value = 144.2
$ws.Range("K23").Value = 432.6
$ws.Range("M23").Value = -197.6

# CUL row 80
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 23943.666
$ws.Range("I80").Value = 29999
$ws.Range("J80").Value = 23186.75
$ws.Range("K80").Value = 89997
$ws.Range("L80").Value = 69560.25
$ws.Range("M80").Value = -89061
$ws.Range("N80").Value = -71432.25

# CUL row 83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 23943.666
$ws.Range("I83").Value = 29999
$ws.Range("J83").Value = 23186.75
$ws.Range("K83").Value = 269991
$ws.Range("L83").Value = 208680.75
$ws.Range("M83").Value = -265311
$ws.Range("N83").Value = -218040.75

# CUL row 86
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 671.1667
$ws.Range("J86").Value = 488.5
$ws.Range("L86").Value = 1465.5
$ws.Range("N86").Value = -3837.5

# CUL row 89
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 671.1667
$ws.Range("J89").Value = 488.5
$ws.Range("L89").Value = 4396.5
$ws.Range("N89").Value = -16252.5

# CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 4606.222
$ws.Range("I92").Value = 184.5
$ws.Range("J92").Value = 13449.667
$ws.Range("K92").Value = 553.5
$ws.Range("L92").Value = 40349.001
$ws.Range("M92").Value = 694.5
$ws.Range("N92").Value = -42845.001

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 6642.4707
$ws.Range("J107").Value = 13541.75
$ws.Range("L107").Value = 40625.25
$ws.Range("N107").Value = -44465.25

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1477.13
$ws.Range("I131").Value = 1115.1428
$ws.Range("J131").Value = 1504.3763
$ws.Range("K131").Value = 3345.4284
$ws.Range("L131").Value = 4513.1289
$ws.Range("M131").Value = 1694.5716
$ws.Range("N131").Value = -14593.1289

# GSM row 52
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 80000
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").Value = $null

# GSM row 58
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 21642.5
$ws.Range("I58").Value = 21642.5
$ws.Range("K58").Value = 21642.5
$ws.Range("M58").Value = -21365.5

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14275.889
$ws.Range("I70").Value = 14995.333
$ws.Range("K70").Value = 14995.333
$ws.Range("M70").Value = -14725.333

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 14275.889
$ws.Range("I73").Value = 14995.333
$ws.Range("K73").Value = 14995.333
$ws.Range("M73").Value = -14059.333

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2400
$ws.Range("I80").Value = 2100
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 2100
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -1102
$ws.Range("N80").Value = -4996

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2400
$ws.Range("I83").Value = 2100
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 10500
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -5508
$ws.Range("N83").Value = -24984

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 740.9583
$ws.Range("I97").Value = 740.9583
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 740.9583
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = $null
$ws.Range("N97").Value = -244.9583

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6606.5
$ws.Range("I7").Value = 7234.9414
$ws.Range("K7").Value = 7234.9414
$ws.Range("M7").Value = -7122.9414

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5332.6665
$ws.Range("I40").Value = 4999.5
$ws.Range("K40").Value = 4999.5
$ws.Range("M40").Value = -4863.5

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3164.0588
$ws.Range("I100").Value = 2677.889
$ws.Range("J100").Value = 3711
$ws.Range("K100").Value = 2677.889
$ws.Range("L100").Value = 3711
$ws.Range("M100").Value = -2136.889
$ws.Range("N100").Value = -4793

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 12219.4
$ws.Range("I122").Value = 13774.25
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 41322.75
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -38872.75
$ws.Range("N122").Value = -22900

# LTW row 125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H125").Value = 45715
$ws.Range("J125").Value = 45715
$ws.Range("L125").Value = 45715
$ws.Range("N125").Value = -55555

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 6606.5
$ws.Range("I126").Value = 7234.9414
$ws.Range("K126").Value = 21704.8242
$ws.Range("M126").Value = -19234.8242

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2136196.2
$ws.Range("I136").Value = 27288.875
$ws.Range("J136").Value = 4546376
$ws.Range("K136").Value = 81866.625
$ws.Range("L136").Value = 13639128
$ws.Range("M136").Value = -79316.625
$ws.Range("N136").Value = -13644228

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4500
$ws.Range("I122").Value = 4500
$ws.Range("K122").Value = 13500
$ws.Range("M122").Value = -11050

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4585.407
$ws.Range("I126").Value = 3428.8096
$ws.Range("J126").Value = 8633.5
$ws.Range("K126").Value = 10286.4288
$ws.Range("L126").Value = 25900.5
$ws.Range("M126").Value = -7816.4288
$ws.Range("N126").Value = -30840.5

# WVR row 127
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 25000
$ws.Range("J127").Value = 25000
$ws.Range("L127").Value = 25000
$ws.Range("N127").Value = -34920

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 710028.2
$ws.Range("I136").Value = 3124.75
$ws.Range("J136").Value = 967084
$ws.Range("K136").Value = 9374.25
$ws.Range("L136").Value = 2901252
$ws.Range("M136").Value = -6824.25
$ws.Range("N136").Value = -2906352
